$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.285.91"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "1.622.99"
$ws.Range("E3").Value = "  +1.45%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("E9").Value = "  +0.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.62%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "1.848.64"
$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("D13").Value = "1.620.77"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("E15").Value = "  +1.05%  "

$ws.Range("D16").Value = "26.296.58"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.75%  "

$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("E22").Value = "  +1.04%  "

$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0527"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.56%  "

$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("E32").Value = "  +2.78%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("E34").Value = "  +1.75%  "

$ws.Range("E35").Value = "  +2.30%  "

$ws.Range("D36").Value = "1.173.13"
$ws.Range("E36").Value = "  +4.39%  "

$ws.Range("E37").Value = "  +1.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.811"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.78%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +1.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.793"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.42%  "

$ws.Range("D44").Value = "1.759.32"
$ws.Range("E44").Value = "  +1.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("D46").Value = "0.0₆0106"
$ws.Range("E46").Value = "  +15.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.410"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("E51").Value = "  -0.14%  "
